$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data down
$ws.Rows.Item(2).Insert()

# Fill in the new row 2 with the new event entry
$ws.Range("A2").Value = "Organizer and moderator"
$ws.Range("B2").Value = "Mar. 5, 2025"
$ws.Range("C2").Value = "International Conference on the Prevention of Sexual Abuse Perpetration: Public Health Perspectives and Challenges"
$ws.Range("D2").Value = "Universidad El Bosque \& Fractales "
$ws.Range("E2").Value = "\href{https://www.youtube.com/watch?v=dNqxY_fGKwE}{Universidad El Bosque}"

$ws.Range("C9").Select() | Out-Null
